$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.428.98"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "'3.589.57"
$ws.Range("E3").Value = "  +4.58%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'238.09"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").Value = "'650.91"
$ws.Range("E6").Value = "  +4.56%  "

$ws.Range("E7").Value = "  +3.92%  "

$ws.Range("D8").Value = "'0.404"
$ws.Range("E8").Value = "  +1.28%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'0.997"
$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").Value = "'3.585.54"
$ws.Range("E11").Value = "  +4.45%  "

$ws.Range("D12").Value = "'42.68"
$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "'6.31"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "'4.275.23"
$ws.Range("E15").Value = "  +4.88%  "

$ws.Range("D16").Value = "'95.348.07"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("D18").Value = "'3.588.17"
$ws.Range("E18").Value = "  +4.78%  "

$ws.Range("D19").Value = "'7.93"
$ws.Range("E19").Value = "  -4.35%  "

$ws.Range("E20").Value = "  +8.51%  "

$ws.Range("D21").Value = "'17.90"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("D22").Value = "'3.63"
$ws.Range("E22").Value = "  +6.74%  "

$ws.Range("E23").Value = "  +3.68%  "

$ws.Range("D24").Value = "'508.92"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  +5.25%  "

$ws.Range("D26").Value = "'6.60"
$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("D27").Value = "'96.39"
$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("D28").Value = "'12.68"
$ws.Range("E28").Value = "  +5.23%  "

$ws.Range("D29").Value = "'3.798.07"
$ws.Range("E29").Value = "  +5.26%  "

$ws.Range("D30").Value = "'3.10"
$ws.Range("E30").Value = "  +12.38%  "

$ws.Range("D31").Value = "'11.33"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.140"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("E35").Value = "  +2.49%  "

$ws.Range("D36").Value = "'31.76"
$ws.Range("E36").Value = "  +5.88%  "

$ws.Range("D37").Value = "'0.557"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").Value = "'8.19"
$ws.Range("E38").Value = "  +8.79%  "

$ws.Range("D39").Value = "'570.01"
$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  +5.39%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "'0.917"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'23.74"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("E45").Value = "  -1.14%  "

$ws.Range("D46").Value = "'5.67"
$ws.Range("E46").Value = "  +2.52%  "

$ws.Range("D47").Value = "'33.94"
$ws.Range("E47").Value = "  +33.84%  "

$ws.Range("D48").Value = "'2.24"
$ws.Range("E48").Value = "  +4.70%  "

$ws.Range("D49").Value = "'0.0414"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").Value = "'3.53"
$ws.Range("E50").Value = "  -4.00%  "

$ws.Range("D51").Value = "'53.77"
